$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: BZDPT6650 -> BZ6650
$ws.Range("G2").Value = "BZ6650"

# E2: "25 Feb" (text) -> real date 2025-03-25, formatted as a date (same numeric date style as K2)
$ws.Range("E2").Value = 45741
$ws.Range("K2").Copy()
$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column E sizing (bestFit/custom width matching column K's width)
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# Update the active selection to H6
$ws.Range("H6").Select() | Out-Null
